# Rename plant_name values in column A to their official legal entity
# names:
#   "Sappi Cloquet" -> "SAPPI CLOQUET LLC"
#   "Boise"         -> "BOISE WHITE PAPER LLC"
#   "Westrock"      -> "WestRock MN Corporation"
#
# (lcoh_policy_modeling_input.xlsx update, 8/5)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$map = @{
    "Sappi Cloquet" = "SAPPI CLOQUET LLC"
    "Boise"         = "BOISE WHITE PAPER LLC"
    "Westrock"      = "WestRock MN Corporation"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    if ($null -ne $old) {
        $key = [string]$old
        if ($map.ContainsKey($key)) {
            $cell.ClearFormats()
            $cell.Value2 = $map[$key]
        }
    }
}

# Leave the selection on the last block of renamed rows, matching the
# editor's final cursor position.
$ws.Range("A170:A172").Select()
